$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Valor Mora" amount and "Cant. Periodos" count near the top of the statement
$ws.Range("E11").Value = 1316800
$ws.Range("F13").Value = 4

# 2) Insert a new row before row 19 - this shifts the old row 18 (with its special
#    bottom-border style) down to row 19, and also shifts the old rows 19-24 down
#    (old row 23 -> 24, old row 24 -> 25) opening a fresh row 18 slot.
$ws.Rows.Item(19).Insert()

# 3) The freshly inserted row 19 is blank; give it the same formatting that row 18
#    had (before we overwrite row 18's own formatting in the next step).
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)  # xlPasteFormats

# 4) Re-format row 18 to match the "normal" data rows (16/17) since it is no longer
#    the last/bottom row of the table.
$ws.Range("B16:J16").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# 5) Re-order/renumber the period values of the existing worker rows and populate the
#    brand-new fourth period row (2508).
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73080914"
$ws.Range("D18").Value = "GUILLERMO ENRIQUE BARRIOS GONZALEZ"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 329200
$ws.Range("G18").Value = 8229977

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73080914"
$ws.Range("D19").Value = "GUILLERMO ENRIQUE BARRIOS GONZALEZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 329200
$ws.Range("G19").Value = 8229977

Write-Host "Done"
